$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Student Name"
$ws.Range("B1").Value = "University Name"
$ws.Range("C1").Value = "Curriculum"
$ws.Range("D1").Value = "Academic Year"
$ws.Range("E1").Value = "Placement Year"
$ws.Range("F1").Value = "Course Year"

$ws.Range("F1").Select()
